$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.530.49'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '''2.493.10'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''570.27'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '''166.49'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").Value = '''4.88'
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").Value = '''2.948.98'
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").Value = '''69.345.69'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").Value = '''24.23'
$ws.Range("E16").Value = '  -2.70%  '
$ws.Range("D17").Value = '''2.501.04'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("D19").Value = '''355.48'
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '''1.90'
$ws.Range("E22").Value = '  -5.24%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '''69.47'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").Value = '''3.82'
$ws.Range("E25").Value = '  -3.75%  '
$ws.Range("D27").Value = '''8.61'
$ws.Range("E27").Value = '  -3.46%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = '''0.0₃0874'
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("E30").Value = '  -2.77%  '
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("D32").Value = '''437.70'
$ws.Range("E32").Value = '  -5.62%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").Value = '''154.68'
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("D36").Value = '''0.113'
$ws.Range("E36").Value = '  -3.62%  '
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("D38").Value = '''18.17'
$ws.Range("E38").Value = '  -2.18%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '''0.314'
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = '''4.59'
$ws.Range("E41").Value = '  -2.83%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.58'
$ws.Range("E42").Value = '  -1.96%  '
$ws.Range("B43").Value = 'POPCAT'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range("D43").Value = '''2.39'
$ws.Range("E43").Value = '  +50.77%  '
$ws.Range("E44").Value = '  -4.22%  '
$ws.Range("D46").Value = '''138.56'
$ws.Range("E46").Value = '  -2.40%  '
$ws.Range("D47").Value = '''3.43'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("E48").Value = '  -3.44%  '
$ws.Range("D49").Value = '''0.0724'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = '''0.572'
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("E51").Value = '  -1.02%  '
